$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Insert the "skeleton" of the 15 new paragraphs (signature blocks) right
# after the current last paragraph of the document body. InsertXML is used
# on a collapsed Range positioned at the very end of the document content so
# the existing content is left untouched and the new paragraphs are appended.
# Run-level character styles (w:rStyle inside a run's w:rPr) are applied in a
# second pass below via Range.Style, because InsertXML silently drops
# run-level rStyle references baked directly into the injected XML.
# Paragraph-mark run properties (w:pPr/w:rPr/w:rStyle) are NOT affected by
# that limitation and so are included directly in the skeleton XML.
# ---------------------------------------------------------------------------

$origCount = $d.Paragraphs.Count

$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)

$skeleton = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">IVAN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>McKEE</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/><w:t>Authorised to sign by the Scottish Ministers</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:t>St Andrew’s House,</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>Edinburgh</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>1st October 2024</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>We consent</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="SigSignee"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:t>JEFF SMITH</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>ANNA TURLEY</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:tab/><w:t>Two of the Lords Commissioners of His Majesty’s Treasury</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="LQN2"/><w:ind w:left="0" w:firstLine="0"/></w:pPr><w:r><w:t>1st October 2024</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($skeleton)

# ---------------------------------------------------------------------------
# Second pass: apply run-level character styles (w:rStyle) on the text that
# needs them. Each new paragraph is located by its index (origCount + N) so
# this does not depend on absolute character offsets into the story.
# ---------------------------------------------------------------------------

function Set-ParaCharStyle($paraIndex, $styleName, $skipStart, $skipEnd) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $start = $r.Start + $skipStart
    $finish = $r.End - $skipEnd
    $sub = $d.Range($start, $finish)
    $sub.Style = $styleName
}

# Paragraph 4 (origCount+4): tab, "IVAN McKEE" -> SigSignee (skip leading tab,
# skip trailing paragraph mark)
Set-ParaCharStyle ($origCount + 4) "SigSignee" 1 1

# Paragraph 5 (origCount+5): tab + "Authorised to sign by the Scottish
# Ministers" -> Sigtitle (whole run incl. tab, excl. paragraph mark)
Set-ParaCharStyle ($origCount + 5) "Sigtitle" 0 1

# Paragraph 6 (origCount+6): "St Andrew's House," -> SigAdd
Set-ParaCharStyle ($origCount + 6) "SigAdd" 0 1

# Paragraph 7 (origCount+7): "Edinburgh" -> SigAdd
Set-ParaCharStyle ($origCount + 7) "SigAdd" 0 1

# Paragraph 8 (origCount+8): "1st October 2024" -> SigDate
Set-ParaCharStyle ($origCount + 8) "SigDate" 0 1

# Paragraph 10 (origCount+10): "We consent" -> Sigsignatory
Set-ParaCharStyle ($origCount + 10) "Sigsignatory" 0 1

# Paragraph 12 (origCount+12): tab, "JEFF SMITH" -> SigSignee (skip leading tab)
Set-ParaCharStyle ($origCount + 12) "SigSignee" 1 1

# Paragraph 13 (origCount+13): tab, "ANNA TURLEY" -> SigSignee (skip leading tab)
Set-ParaCharStyle ($origCount + 13) "SigSignee" 1 1

# Paragraph 14 (origCount+14): tab + "Two of the Lords Commissioners of His
# Majesty's Treasury" -> Sigtitle (whole run incl. tab)
Set-ParaCharStyle ($origCount + 14) "Sigtitle" 0 1

# Paragraph 15 (origCount+15): "1st October 2024" -> SigDate
Set-ParaCharStyle ($origCount + 15) "SigDate" 0 1

Write-Host "New paragraph count: $($d.Paragraphs.Count)"
